$wb = $excel.ActiveWorkbook

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1635.3334
$ws.Range("I43").Value = 942.8570999999999
$ws.Range("J43").Value = 1920.4706
$ws.Range("K43").Value = 942.8570999999999
$ws.Range("L43").Value = 1920.4706
$ws.Range("M43").Value = -873.8570999999999
$ws.Range("N43").Value = -2058.4706
$ws.Range("H74").Value = 3963.3333
$ws.Range("I74").Value = 3945
$ws.Range("K74").Value = 3945
$ws.Range("M74").Value = -3009
$ws.Range("H77").Value = 3963.3333
$ws.Range("I77").Value = 3945
$ws.Range("K77").Value = 19725
$ws.Range("M77").Value = -15045
$ws.Range("H132").Value = 804664.7
$ws.Range("I132").Value = 1172.2452
$ws.Range("J132").Value = 6127802
$ws.Range("K132").Value = 3516.7356
$ws.Range("L132").Value = 18383406
$ws.Range("M132").Value = -986.7356
$ws.Range("N132").Value = -18388466
$ws.Range("H137").Value = 3849271.2
$ws.Range("I137").Value = 4547738.5
$ws.Range("K137").Value = 13643215.5
$ws.Range("M137").Value = -13640665.5
$ws.Range("H141").Value = 1808.15
$ws.Range("I141").Value = 885.1875
$ws.Range("J141").Value = 5500
$ws.Range("K141").Value = 2655.5625
$ws.Range("L141").Value = 16500
$ws.Range("M141").Value = 2524.4375
$ws.Range("N141").Value = -26860

# ARM sheet updates
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1978.8334
$ws.Range("I2").Value = 1876.4117
$ws.Range("J2").Value = 2227.5715
$ws.Range("K2").Value = 1876.4117
$ws.Range("L2").Value = 2227.5715
$ws.Range("M2").Value = -1763.4117
$ws.Range("N2").Value = -2453.5715
$ws.Range("H16").Value = 10001.667
$ws.Range("I16").Value = 10001.667
$ws.Range("K16").Value = 10001.667
$ws.Range("M16").Value = -9714.666999999999
$ws.Range("H32").Value = 15202.255
$ws.Range("I32").Value = 18314.955
$ws.Range("J32").Value = 7061.346
$ws.Range("K32").Value = 18314.955
$ws.Range("L32").Value = 7061.346
$ws.Range("M32").Value = -18027.955
$ws.Range("N32").Value = -7635.346
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10626
$ws.Range("H74").Value = 3707407.5
$ws.Range("I74").Value = 4483106
$ws.Range("J74").Value = 87480.5
$ws.Range("K74").Value = 4483106
$ws.Range("L74").Value = 87480.5
$ws.Range("M74").Value = -4482232
$ws.Range("N74").Value = -89228.5
$ws.Range("H77").Value = 3707407.5
$ws.Range("I77").Value = 4483106
$ws.Range("J77").Value = 87480.5
$ws.Range("K77").Value = 22415530
$ws.Range("L77").Value = 437402.5
$ws.Range("M77").Value = -22411162
$ws.Range("N77").Value = -446138.5
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H116").Value = 1978.8334
$ws.Range("I116").Value = 1876.4117
$ws.Range("J116").Value = 2227.5715
$ws.Range("K116").Value = 1876.4117
$ws.Range("L116").Value = 2227.5715
$ws.Range("M116").Value = 417.5882999999999
$ws.Range("N116").Value = -6815.5715

# BSM sheet updates
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1978.8334
$ws.Range("I3").Value = 1876.4117
$ws.Range("J3").Value = 2227.5715
$ws.Range("K3").Value = 1876.4117
$ws.Range("L3").Value = 2227.5715
$ws.Range("M3").Value = -1762.4117
$ws.Range("N3").Value = -2455.5715
$ws.Range("H20").Value = 1248.6364
$ws.Range("I20").Value = 983.2
$ws.Range("J20").Value = 1469.8334
$ws.Range("K20").Value = 983.2
$ws.Range("L20").Value = 1469.8334
$ws.Range("M20").Value = -736.2
$ws.Range("N20").Value = -1963.8334
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 1200
$ws.Range("I94").Value = 1400
$ws.Range("K94").Value = 1400
$ws.Range("M94").Value = -949
$ws.Range("H130").Value = 38998
$ws.Range("J130").Value = 38998
$ws.Range("L130").Value = 38998
$ws.Range("N130").Value = -49038
$ws.Range("H137").Value = 70000
$ws.Range("J137").Value = 70000
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200

# CRP sheet updates
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 19447.5
$ws.Range("I12").Value = 930
$ws.Range("K12").Value = 930
$ws.Range("M12").Value = -760
$ws.Range("H31").Value = 2478.28
$ws.Range("I31").Value = 1430.6666
$ws.Range("J31").Value = 6192.5454
$ws.Range("K31").Value = 1430.6666
$ws.Range("L31").Value = 6192.5454
$ws.Range("M31").Value = -1135.6666
$ws.Range("N31").Value = -6782.5454
$ws.Range("H34").Value = 2478.28
$ws.Range("I34").Value = 1430.6666
$ws.Range("J34").Value = 6192.5454
$ws.Range("K34").Value = 1430.6666
$ws.Range("L34").Value = 6192.5454
$ws.Range("M34").Value = -1228.6666
$ws.Range("N34").Value = -6596.5454

# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 156.22223
$ws.Range("I38").Value = 33.333332
$ws.Range("J38").Value = 217.66667
$ws.Range("K38").Value = 99.999996
$ws.Range("L38").Value = 653.00001
$ws.Range("M38").Value = 247.000004
$ws.Range("N38").Value = -1347.00001
$ws.Range("H131").Value = 1345.4615
$ws.Range("J131").Value = 1359.8158
$ws.Range("L131").Value = 4079.4474
$ws.Range("N131").Value = -14159.4474

# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 25464.428
$ws.Range("J93").Value = 30083.666
$ws.Range("L93").Value = 30083.666
$ws.Range("N93").Value = -33827.666
$ws.Range("H97").Value = 1394.4706
$ws.Range("I97").Value = 1465.1333
$ws.Range("J97").Value = 864.5
$ws.Range("K97").Value = 1465.1333
$ws.Range("L97").Value = 864.5
$ws.Range("M97").Value = -969.1333
$ws.Range("N97").Value = -1856.5
$ws.Range("H122").Value = 2731.0286
$ws.Range("I122").Value = 2354.4443
$ws.Range("J122").Value = 4002
$ws.Range("K122").Value = 7063.3329
$ws.Range("L122").Value = 12006
$ws.Range("M122").Value = -4613.3329
$ws.Range("N122").Value = -16906
$ws.Range("H136").Value = 17681.5
$ws.Range("J136").Value = 17681.5
$ws.Range("L136").Value = 53044.5
$ws.Range("N136").Value = -58144.5

# LTW sheet updates
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2666.6667
$ws.Range("I93").Value = 2500
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2500
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1252
$ws.Range("N93").Value = -5496
$ws.Range("H122").Value = 3344.5715
$ws.Range("I122").Value = 2988.5
$ws.Range("J122").Value = 3487
$ws.Range("K122").Value = 8965.5
$ws.Range("L122").Value = 10461
$ws.Range("M122").Value = -6515.5
$ws.Range("N122").Value = -15361
